# Confusion matrices for 8s window
# - Clear the stray confusion-matrix summary text previously dumped in I5:I8
#   (and the now-unused shared strings that backed them)
# - Update the selection/scroll position left over from editing
# - Re-anchor the shared formula in F22:G22 (was incorrectly spanning B22:G22)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8s window")

$ws.Range("I5").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("I8").ClearContents()

$ws.Range("F22").Formula = "=SUM(F2:F21)/COUNT(F2:F21)"
$ws.Range("G22").Formula = "=SUM(G2:G21)/COUNT(G2:G21)"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E13").Select()
